$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 746.3889
$ws.Range("J38").Value = 1562.5
$ws.Range("L38").Value = 4687.5
$ws.Range("N38").Value = -5431.5
$ws.Range("H43").Value = 1154.1818
$ws.Range("J43").Value = 1070.375
$ws.Range("L43").Value = 1070.375
$ws.Range("N43").Value = -1208.375
$ws.Range("H62").Value = 181822290
$ws.Range("I62").Value = 111114690
$ws.Range("K62").Value = 111114690
$ws.Range("M62").Value = -111114066
$ws.Range("H65").Value = 181822290
$ws.Range("I65").Value = 111114690
$ws.Range("K65").Value = 555573450
$ws.Range("M65").Value = -555570330
$ws.Range("H74").Value = 4172.5884
$ws.Range("I74").Value = 4404.8887
$ws.Range("J74").Value = 3911.25
$ws.Range("K74").Value = 4404.8887
$ws.Range("L74").Value = 3911.25
$ws.Range("M74").Value = -3468.8887
$ws.Range("N74").Value = -5783.25
$ws.Range("H77").Value = 4172.5884
$ws.Range("I77").Value = 4404.8887
$ws.Range("J77").Value = 3911.25
$ws.Range("K77").Value = 22024.4435
$ws.Range("L77").Value = 19556.25
$ws.Range("M77").Value = -17344.4435
$ws.Range("N77").Value = -28916.25
$ws.Range("H88").Value = 4907.7334
$ws.Range("I88").Value = 5148.727
$ws.Range("J88").Value = 4245
$ws.Range("K88").Value = 5148.727
$ws.Range("L88").Value = 4245
$ws.Range("M88").Value = -4742.727
$ws.Range("N88").Value = -5057
$ws.Range("H91").Value = 4907.7334
$ws.Range("I91").Value = 5148.727
$ws.Range("J91").Value = 4245
$ws.Range("K91").Value = 5148.727
$ws.Range("L91").Value = 4245
$ws.Range("M91").Value = -3744.727
$ws.Range("N91").Value = -7053
$ws.Range("H111").Value = 2474
$ws.Range("I111").Value = 3005.6
$ws.Range("J111").Value = 1588
$ws.Range("K111").Value = 9016.799999999999
$ws.Range("L111").Value = 4764
$ws.Range("M111").Value = -5949.799999999999
$ws.Range("N111").Value = -10898
$ws.Range("H113").Value = 2909.7896
$ws.Range("I113").Value = 2651.25
$ws.Range("J113").Value = 2978.7334
$ws.Range("K113").Value = 2651.25
$ws.Range("L113").Value = 2978.7334
$ws.Range("M113").Value = 602.75
$ws.Range("N113").Value = -9486.733400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2149.9285
$ws.Range("I45").Value = 2184.7693
$ws.Range("J45").Value = 1697
$ws.Range("K45").Value = 2184.7693
$ws.Range("L45").Value = 1697
$ws.Range("M45").Value = -1807.7693
$ws.Range("N45").Value = -2451
$ws.Range("H122").Value = 1618.5483
$ws.Range("I122").Value = 1346.7037
$ws.Range("J122").Value = 3453.5
$ws.Range("K122").Value = 4040.1111
$ws.Range("L122").Value = 10360.5
$ws.Range("M122").Value = -1590.1111
$ws.Range("N122").Value = -15260.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1212.2941
$ws.Range("I99").Value = 1116.0769
$ws.Range("K99").Value = 1116.0769
$ws.Range("M99").Value = 381.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 496.5
$ws.Range("I22").Value = 443.33334
$ws.Range("J22").Value = 540
$ws.Range("K22").Value = 443.33334
$ws.Range("L22").Value = 540
$ws.Range("M22").Value = -93.33334000000002
$ws.Range("N22").Value = -1240
$ws.Range("H31").Value = 6915132
$ws.Range("I31").Value = 5584277
$ws.Range("J31").Value = 9092895
$ws.Range("K31").Value = 5584277
$ws.Range("L31").Value = 9092895
$ws.Range("M31").Value = -5583982
$ws.Range("N31").Value = -9093485
$ws.Range("H34").Value = 6915132
$ws.Range("I34").Value = 5584277
$ws.Range("J34").Value = 9092895
$ws.Range("K34").Value = 5584277
$ws.Range("L34").Value = 9092895
$ws.Range("M34").Value = -5584075
$ws.Range("N34").Value = -9093299
$ws.Range("H99").Value = 8206.429
$ws.Range("I99").Value = 8489
$ws.Range("K99").Value = 8489
$ws.Range("M99").Value = -6991
$ws.Range("H126").Value = 8206.429
$ws.Range("I126").Value = 8489
$ws.Range("K126").Value = 25467
$ws.Range("M126").Value = -22997
$ws.Range("H134").Value = 3146.7334
$ws.Range("I134").Value = 3032.36
$ws.Range("J134").Value = 3718.6
$ws.Range("K134").Value = 9097.08
$ws.Range("L134").Value = 11155.8
$ws.Range("M134").Value = -6562.08
$ws.Range("N134").Value = -16225.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 494.64517
$ws.Range("I5").Value = 516.62964
$ws.Range("J5").Value = 346.25
$ws.Range("K5").Value = 1549.88892
$ws.Range("L5").Value = 1038.75
$ws.Range("M5").Value = -1437.88892
$ws.Range("N5").Value = -1262.75
$ws.Range("H101").Value = 16144.444
$ws.Range("J101").Value = 16144.444
$ws.Range("L101").Value = 48433.33199999999
$ws.Range("N101").Value = -53301.33199999999
$ws.Range("H122").Value = 498.6
$ws.Range("I122").Value = 372
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 3348
$ws.Range("L122").Value = 9045
$ws.Range("M122").Value = -898
$ws.Range("N122").Value = -13945
$ws.Range("H135").Value = 494.64517
$ws.Range("I135").Value = 516.62964
$ws.Range("J135").Value = 346.25
$ws.Range("K135").Value = 4649.66676
$ws.Range("L135").Value = 3116.25
$ws.Range("M135").Value = -2114.66676
$ws.Range("N135").Value = -8186.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 1221
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1221
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 1221
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -1797
$ws.Range("H34").Value = 35668.25
$ws.Range("J34").Value = 35668.25
$ws.Range("L34").Value = 35668.25
$ws.Range("N34").Value = -36204.25
$ws.Range("H76").Value = 35668.25
$ws.Range("J76").Value = 35668.25
$ws.Range("L76").Value = 35668.25
$ws.Range("N76").Value = -36298.25
$ws.Range("H79").Value = 35668.25
$ws.Range("J79").Value = 35668.25
$ws.Range("L79").Value = 35668.25
$ws.Range("N79").Value = -37852.25
$ws.Range("H80").Value = 2747.625
$ws.Range("J80").Value = 2788.9092
$ws.Range("L80").Value = 2788.9092
$ws.Range("N80").Value = -4784.9092
$ws.Range("H81").Value = 1221
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 1221
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 1221
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -3217
$ws.Range("H83").Value = 2747.625
$ws.Range("J83").Value = 2788.9092
$ws.Range("L83").Value = 13944.546
$ws.Range("N83").Value = -23928.546
$ws.Range("H84").Value = 1221
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 1221
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 3663
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -13647
$ws.Range("H102").Value = 2819.276
$ws.Range("I102").Value = 3067.7827
$ws.Range("J102").Value = 1866.6666
$ws.Range("K102").Value = 3067.7827
$ws.Range("L102").Value = 1866.6666
$ws.Range("M102").Value = -1445.7827
$ws.Range("N102").Value = -5110.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2310.5
$ws.Range("I40").Value = 1998
$ws.Range("J40").Value = 2373
$ws.Range("K40").Value = 1998
$ws.Range("L40").Value = 2373
$ws.Range("M40").Value = -1862
$ws.Range("N40").Value = -2645
$ws.Range("H46").Value = 633.3333
$ws.Range("J46").Value = 700
$ws.Range("L46").Value = 700
$ws.Range("N46").Value = -1076
$ws.Range("H61").Value = 1645.2727
$ws.Range("I61").Value = 1559.8
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 1559.8
$ws.Range("L61").Value = 2500
$ws.Range("M61").Value = -1357.8
$ws.Range("N61").Value = -2904
$ws.Range("H113").Value = 1645.2727
$ws.Range("I113").Value = 1559.8
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1559.8
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 610.2
$ws.Range("N113").Value = -6840
$ws.Range("H122").Value = 4250.6665
$ws.Range("I122").Value = 4699.8
$ws.Range("J122").Value = 2005
$ws.Range("K122").Value = 14099.4
$ws.Range("L122").Value = 6015
$ws.Range("M122").Value = -11649.4
$ws.Range("N122").Value = -10915

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 34484804
$ws.Range("I81").Value = 58825236
$ws.Range("J81").Value = 2523.75
$ws.Range("K81").Value = 117650472
$ws.Range("L81").Value = 5047.5
$ws.Range("M81").Value = -117649411
$ws.Range("N81").Value = -7169.5
$ws.Range("H84").Value = 34484804
$ws.Range("I84").Value = 58825236
$ws.Range("J84").Value = 2523.75
$ws.Range("K84").Value = 588252360
$ws.Range("L84").Value = 25237.5
$ws.Range("M84").Value = -588247056
$ws.Range("N84").Value = -35845.5
$ws.Range("H122").Value = 83335230
$ws.Range("I122").Value = 125001624
$ws.Range("J122").Value = 2438.75
$ws.Range("K122").Value = 375004872
$ws.Range("L122").Value = 2438.75
$ws.Range("M122").Value = -375002422
$ws.Range("N122").Value = -12216.25
